$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six oldest year columns (B:G = 2009..2014). This shifts the
# existing H:M (2015..2020) data left into B:G, matching the diff which
# drops 2009-2014 and keeps 2015-2020.
$ws.Range("B:G").Delete() | Out-Null

# A few of the recomputed statistics in the final (2020) column differ
# slightly from the previous run - update them to the new values.
$ws.Range("G3").Value = 1.851116598963962
$ws.Range("G4").Value = 1.223013465149995
$ws.Range("G6").Value = 0.968754992699923
